# Updates leveling-profit market data (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the eight job sheets, driven by a scheduled market-data refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70: Consecrating Congregation / Holy Water
$ws.Range("H70").Value = 1468.0769
$ws.Range("J70").Value = 1550
$ws.Range("L70").Value = 4650
$ws.Range("N70").Value = -5190

# Row 73: Curbing the Contagion (L) / Holy Water
$ws.Range("H73").Value = 1468.0769
$ws.Range("J73").Value = 1550
$ws.Range("L73").Value = 4650
$ws.Range("N73").Value = -6522

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 1050.2354
$ws.Range("I112").Value = 690
$ws.Range("J112").Value = 1112.3448
$ws.Range("K112").Value = 2070
$ws.Range("L112").Value = 3337.0344
$ws.Range("M112").Value = -962
$ws.Range("N112").Value = -5553.0344

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 4363.273
$ws.Range("I116").Value = 2533.3333
$ws.Range("J116").Value = 5049.5
$ws.Range("K116").Value = 2533.3333
$ws.Range("L116").Value = 5049.5
$ws.Range("M116").Value = 908.6667000000002
$ws.Range("N116").Value = -11933.5

# Row 125: Body over Mind / Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 238.66667
$ws.Range("I125").Value = 265.25
$ws.Range("J125").Value = 208.28572
$ws.Range("K125").Value = 2387.25
$ws.Range("L125").Value = 1874.57148
$ws.Range("M125").Value = 72.75
$ws.Range("N125").Value = -6794.571480000001

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 2107.0188
$ws.Range("I132").Value = 2160.4082
$ws.Range("J132").Value = 1453
$ws.Range("K132").Value = 6481.2246
$ws.Range("L132").Value = 4359
$ws.Range("M132").Value = -3951.2246
$ws.Range("N132").Value = -9419

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 611.64514
$ws.Range("I2").Value = 703
$ws.Range("J2").Value = 419.8
$ws.Range("K2").Value = 703
$ws.Range("L2").Value = 419.8
$ws.Range("M2").Value = -590
$ws.Range("N2").Value = -645.8

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 8756.313
$ws.Range("I32").Value = 6279.4634
$ws.Range("J32").Value = 18911.4
$ws.Range("K32").Value = 6279.4634
$ws.Range("L32").Value = 18911.4
$ws.Range("M32").Value = -5992.4634
$ws.Range("N32").Value = -19485.4

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 2484.6365
$ws.Range("I45").Value = 1714.1428
$ws.Range("J45").Value = 3833
$ws.Range("K45").Value = 1714.1428
$ws.Range("L45").Value = 3833
$ws.Range("M45").Value = -1337.1428
$ws.Range("N45").Value = -4587

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 611.64514
$ws.Range("I116").Value = 703
$ws.Range("J116").Value = 419.8
$ws.Range("K116").Value = 703
$ws.Range("L116").Value = 419.8
$ws.Range("M116").Value = 1591
$ws.Range("N116").Value = -5007.8

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 26985.95
$ws.Range("I132").Value = 2602.1
$ws.Range("J132").Value = 51369.8
$ws.Range("K132").Value = 7806.299999999999
$ws.Range("L132").Value = 154109.4
$ws.Range("M132").Value = -5276.299999999999
$ws.Range("N132").Value = -159169.4

# Row 139: Backing up My Words / Titanium Gold Thornplate of Fending
$ws.Range("H139").Value = 40264.9
$ws.Range("J139").Value = 40264.9
$ws.Range("L139").Value = 40264.9
$ws.Range("N139").Value = -50544.9

# Row 140: A Hand for a Deckhand / Ra'Kaznar Gloves of Scouting
$ws.Range("H140").Value = 53445.6
$ws.Range("J140").Value = 53445.6
$ws.Range("L140").Value = 53445.6
$ws.Range("N140").Value = -63805.6

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 611.64514
$ws.Range("I3").Value = 703
$ws.Range("J3").Value = 419.8
$ws.Range("K3").Value = 703
$ws.Range("L3").Value = 419.8
$ws.Range("M3").Value = -589
$ws.Range("N3").Value = -647.8

# Row 81: Diamond Sawdust / Titanium Battleaxe
$ws.Range("H81").Value = 10635.454
$ws.Range("J81").Value = 10635.454
$ws.Range("L81").Value = 10635.454
$ws.Range("N81").Value = -12757.454

# Row 84: I'm a Lumberjack and I'm Okay (L) / Titanium Battleaxe
$ws.Range("H84").Value = 10635.454
$ws.Range("J84").Value = 10635.454
$ws.Range("L84").Value = 31906.362
$ws.Range("N84").Value = -42514.362

# Row 100: And My Axe / Doman Iron War Axe
$ws.Range("H100").Value = 13499.5
$ws.Range("J100").Value = 13499.5
$ws.Range("L100").Value = 13499.5
$ws.Range("N100").Value = -15663.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 4167.8335
$ws.Range("I31").Value = 4202.8887
$ws.Range("J31").Value = 4152.8096
$ws.Range("K31").Value = 4202.8887
$ws.Range("L31").Value = 4152.8096
$ws.Range("M31").Value = -3907.8887
$ws.Range("N31").Value = -4742.8096

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 4167.8335
$ws.Range("I34").Value = 4202.8887
$ws.Range("J34").Value = 4152.8096
$ws.Range("K34").Value = 4202.8887
$ws.Range("L34").Value = 4152.8096
$ws.Range("M34").Value = -4000.8887
$ws.Range("N34").Value = -4556.8096

# Row 95: Standing on Ceremony / High Steel Fork
$ws.Range("H95").Value = 28333.334
$ws.Range("J95").Value = 28333.334
$ws.Range("L95").Value = 28333.334
$ws.Range("N95").Value = -33825.334

# Row 96: Composition / Larch Composite Bow
$ws.Range("H96").Value = 17081
$ws.Range("J96").Value = 17081
$ws.Range("L96").Value = 17081
$ws.Range("N96").Value = -22573

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 20836584
$ws.Range("I99").Value = 2722.2222
$ws.Range("K99").Value = 2722.2222
$ws.Range("M99").Value = -1224.2222

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 20836584
$ws.Range("I126").Value = 2722.2222
$ws.Range("K126").Value = 8166.6666
$ws.Range("M126").Value = -5696.6666

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 2190.4443
$ws.Range("I5").Value = 1218
$ws.Range("J5").Value = 2676.6667
$ws.Range("K5").Value = 3654
$ws.Range("L5").Value = 8030.000100000001
$ws.Range("M5").Value = -3542
$ws.Range("N5").Value = -8254.000100000001

# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 11878
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 11878
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 35634
$ws.Range("N68").Value = -37256
$ws.Range("M68").ClearContents()

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 11878
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 11878
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 106902
$ws.Range("N71").Value = -115014
$ws.Range("M71").ClearContents()

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 737.52
$ws.Range("I131").Value = 300
$ws.Range("J131").Value = 741.9394
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 2225.8182
$ws.Range("M131").Value = 4140
$ws.Range("N131").Value = -12305.8182

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 2190.4443
$ws.Range("I135").Value = 1218
$ws.Range("J135").Value = 2676.6667
$ws.Range("K135").Value = 10962
$ws.Range("L135").Value = 24090.0003
$ws.Range("M135").Value = -8427
$ws.Range("N135").Value = -29160.0003

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 4179.9473
$ws.Range("I80").Value = 3303.3333
$ws.Range("K80").Value = 3303.3333
$ws.Range("M80").Value = -2305.3333

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 4179.9473
$ws.Range("I83").Value = 3303.3333
$ws.Range("K83").Value = 16516.6665
$ws.Range("M83").Value = -11524.6665

# Row 92: Play It by Ear / Triphane Earrings of Healing
$ws.Range("H92").Value = 13785.857
$ws.Range("J92").Value = 13785.857
$ws.Range("L92").Value = 13785.857
$ws.Range("N92").Value = -17529.857

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 2860.2104
$ws.Range("I7").Value = 2449.7058
$ws.Range("K7").Value = 2449.7058
$ws.Range("M7").Value = -2337.7058

# Row 97: Looking for Glove / Gyuki Leather Halfgloves of Scouting
$ws.Range("H97").Value = 14672
$ws.Range("J97").Value = 14672
$ws.Range("L97").Value = 14672
$ws.Range("N97").Value = -16654

# Row 104: Brace Yourselves / Gazelleskin Bracers of Fending
$ws.Range("H104").Value = 24539.8
$ws.Range("J104").Value = 24539.8
$ws.Range("L104").Value = 24539.8
$ws.Range("N104").Value = -31527.8

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 2860.2104
$ws.Range("I126").Value = 2449.7058
$ws.Range("K126").Value = 7349.117400000001
$ws.Range("M126").Value = -4879.117400000001

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 525572.2
$ws.Range("I132").Value = 861627.1
$ws.Range("J132").Value = 2820
$ws.Range("K132").Value = 2584881.3
$ws.Range("L132").Value = 8460
$ws.Range("M132").Value = -2582351.3
$ws.Range("N132").Value = -13520

$ws = $wb.Worksheets.Item("WVR")
# Row 69: Fashion Patrol / Holy Rainbow Sarouel of Casting
$ws.Range("H69").Value = 16817.75
$ws.Range("J69").Value = 16817.75
$ws.Range("L69").Value = 16817.75
$ws.Range("N69").Value = -18315.75

# Row 72: Dress Code Violation (L) / Holy Rainbow Sarouel of Casting
$ws.Range("H72").Value = 16817.75
$ws.Range("J72").Value = 16817.75
$ws.Range("L72").Value = 50453.25
$ws.Range("N72").Value = -57941.25

# Row 101: Who War It Better / Serge Hose of Aiming
$ws.Range("H101").Value = 12676.25
$ws.Range("J101").Value = 12676.25
$ws.Range("L101").Value = 12676.25
$ws.Range("N101").Value = -19166.25

# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 1080.2354
$ws.Range("I122").Value = 1016.0417
$ws.Range("J122").Value = 1234.3
$ws.Range("K122").Value = 3048.1251
$ws.Range("L122").Value = 3702.9
$ws.Range("M122").Value = -598.1251000000002
$ws.Range("N122").Value = -8602.9

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1234.45
$ws.Range("I132").Value = 1000.5
$ws.Range("J132").Value = 2170.25
$ws.Range("K132").Value = 3001.5
$ws.Range("L132").Value = 6510.75
$ws.Range("M132").Value = -471.5
$ws.Range("N132").Value = -11570.75
